$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 39) following the same pattern as the
# preceding rows: date_of_forecast, y_0, y_0_forecast, y_1, y_1_forecast
$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.3398512689293476
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = 0.8571438361188566

# Match the date-column formatting used by the rows above (row 38) so the
# new cell carries the same style (bordered, bold, centered, custom date
# number format) instead of the workbook default style.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
